# [이용섭] Add - [Table] ItemTable, BasePathTable 수정, 아이템 에셋 - 테이블 바인딩
#
# Adds a new "PathFile" column (column H) to the SaleStand table:
#   H2 -> header "PathFile"
#   H3 -> type row value "int32" (matches the other type cells on row 3)
# and leaves the selection on the newly-edited area (J6), matching the
# author's last-saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header + type-row entry for the PathFile field.
$ws.Range("H2").Value = "PathFile"
$ws.Range("H3").Value = "int32"

# Leave the selection where the author left it after the edit.
$ws.Range("J6").Select()
